# New Charger IC, Added Ground Plane
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Include in PCB" column (column J) entirely - no longer tracked.
$ws.Range("J1:J18").EntireColumn.Delete() | Out-Null

# Swap the battery charger IC (row 16 / U5) from the TPOWER TP4056 part to the
# new 3PEAK TPB4056B2X-ES1R part sourced from Digikey.
$ws.Range("D16").Value = "3PEAK"
$ws.Range("E16").Value = "Digikey"
$ws.Range("F16").Value = "TPB4056B2X-ES1R"
$ws.Range("G16").Value = "Charger IC Lithium Ion/Polymer 8-ESOP"
$ws.Range("H16").Value = "ESOP8"

# Update the selected cell to reflect where the author was working (e.g. near
# the newly added ground plane note) when the workbook was saved.
$ws.Range("G25").Select() | Out-Null
